# Remove pi header from BOM.
# Row 2 of the BOM sheet holds the Raspberry Pi Zero WH entry
# (ADA3708 / Raspberry Pi Zero WH _Zero W with Headers_ / A1 / ADA3708_RPI-ZERO / C50982).
# Deleting that row shifts every following row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

# Select the new row 2 (mirrors the sheetView selection captured after the edit).
$ws.Range("A2:XFD2").Select()
